$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. TopicPage: change the remembered selection from a single active
#    cell (A6) to the full data range (A1:B5). Do this before adding /
#    activating any other sheet so TopicPage doesn't end up being the
#    last-active (tabSelected) sheet.
# ------------------------------------------------------------------
$topicPage = $wb.Worksheets.Item("TopicPage")
$topicPage.Range("A1:B5").Select()

# ------------------------------------------------------------------
# 2. Add the new "ErrorPage" worksheet right after "DynamicListingPage"
#    (and before "HomePage"), matching the new sheet order:
#    BlogPostPage, BlogSeriesPage, CTHPPage, DynamicListingPage,
#    ErrorPage, HomePage, InnerPage, LandingPage, TopicPage.
#
#    Copy "HomePage" as a template so the new sheet inherits the same
#    two-column layout, bold/shaded header style and page setup as the
#    rest of the workbook, then becomes the active sheet (as intended -
#    it is the sheet the workbook opens to, activeTab points at it).
# ------------------------------------------------------------------
$dynamicListingPage = $wb.Worksheets.Item("DynamicListingPage")
$homePage = $wb.Worksheets.Item("HomePage")
$homePage.Copy($null, $dynamicListingPage)
$errorPage = $wb.Worksheets.Item("HomePage (2)")
$errorPage.Name = "ErrorPage"

# Drop the copied HomePage rows and add the Error Page test data.
$errorPage.Range("A2:B4").ClearContents()

# Enter values in the same order the original authored the sheet so the
# shared-string table is built up identically.
$errorPage.Range("B3").Value = "Error Page"
$errorPage.Range("B4").Value = "Thank You"
$errorPage.Range("B2").Value = "Page Not Found"
$errorPage.Range("A3").Value = "/PublishedContent/ErrorMessages/Error.html"
$errorPage.Range("A4").Value = "/PublishedContent/ErrorMessages/ThankYou.html"
$errorPage.Range("A2").Value = "/PublishedContent/ErrorMessages/PageNotFound.html"

# Size the columns to fit the new (longer) path / content-type strings.
$errorPage.Columns.Item(1).ColumnWidth = 50.877604166666664
$errorPage.Columns.Item(2).ColumnWidth = 14.307291666666666

$errorPage.Range("A5").Select()
